$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Added user input for indicator"
$ws.Range("B2").Value = "Works with indicator, shows percentage and games won"
$ws.Range("C2").Value = "Ran with inputed indicator and showed results"
$ws.Range("D2").Value = "Works with multiple indicators"

$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Value = Get-Date -Year 2018 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Range("D2").Select()
